$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 141, shifting rows 141:170 down to 142:171
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new weekly record
$ws.Cells.Item(141, 1).Value = 7
$ws.Cells.Item(141, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(141, 3).Value = 'Ñuble'
$ws.Cells.Item(141, 4).Value = 45135
$ws.Cells.Item(141, 5).Value = 16
$ws.Cells.Item(141, 6).Value = 'Fruta'
$ws.Cells.Item(141, 7).Value = 100108
$ws.Cells.Item(141, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(141, 9).Value = 100108002
$ws.Cells.Item(141, 10).Value = 'Mango'
$ws.Cells.Item(141, 11).Value = 'Sin especificar'
$ws.Cells.Item(141, 12).Value = 'Primera'
$ws.Cells.Item(141, 13).Value = 30
$ws.Cells.Item(141, 14).Value = 8000
$ws.Cells.Item(141, 15).Value = 8000
$ws.Cells.Item(141, 16).Value = 8000
$ws.Cells.Item(141, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(141, 18).Value = 'Brasil'
$ws.Cells.Item(141, 19).Value = 2000
$ws.Cells.Item(141, 20).Value = 4

# Ensure the date cell keeps the same date number format as the rest of column D
$ws.Cells.Item(141, 4).NumberFormat = $ws.Cells.Item(142, 4).NumberFormat
